# Adds 5 new strategy worksheets (strategy_id-5011 .. strategy_id-5015) to the workbook,
# cloned from the strategy_id-5009 template sheet (same header row + row-2 structure),
# then overwrites the row-2 trajectory values (columns U:AS) with the per-sheet calibrated values.

$wb = $excel.ActiveWorkbook
$template = $wb.Worksheets.Item("strategy_id-5009")

# --- Add strategy_id-5011 ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "strategy_id-5011"
$template.Range("A1:AS2").Copy($newSheet.Range("A1:AS2"))

$values_5 = @{
    "U2" = 0.9829063213546697
    "V2" = 0.9658126427093395
    "W2" = 0.9487189640640092
    "X2" = 0.931625285418679
    "Y2" = 0.9145316067733488
    "Z2" = 0.8974379281280185
    "AA2" = 0.8803442494826883
    "AB2" = 0.863250570837358
    "AC2" = 0.8461568921920278
    "AD2" = 0.8290632135466975
    "AE2" = 0.8119695349013674
    "AF2" = 0.7948758562560371
    "AG2" = 0.7777821776107068
    "AH2" = 0.7606884989653766
    "AI2" = 0.7435948203200462
    "AJ2" = 0.726501141674716
    "AK2" = 0.7094074630293857
    "AL2" = 0.6923137843840557
    "AM2" = 0.6752201057387253
    "AN2" = 0.6581264270933951
    "AO2" = 0.6410327484480649
    "AP2" = 0.6239390698027346
    "AQ2" = 0.6068453911574043
    "AR2" = 0.5897517125120741
    "AS2" = 0.5726580338667439
}
foreach ($key in $values_5.Keys) {
    $newSheet.Range($key).Value = $values_5[$key]
}

# --- Add strategy_id-5012 ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "strategy_id-5012"
$template.Range("A1:AS2").Copy($newSheet.Range("A1:AS2"))

$values_6 = @{
    "U2" = 0.9862304763364699
    "V2" = 0.9724609526729397
    "W2" = 0.9586914290094095
    "X2" = 0.9449219053458794
    "Y2" = 0.9311523816823493
    "Z2" = 0.9173828580188191
    "AA2" = 0.903613334355289
    "AB2" = 0.8898438106917588
    "AC2" = 0.8760742870282288
    "AD2" = 0.8623047633646985
    "AE2" = 0.8485352397011684
    "AF2" = 0.8347657160376383
    "AG2" = 0.820996192374108
    "AH2" = 0.8072266687105779
    "AI2" = 0.7934571450470478
    "AJ2" = 0.7796876213835175
    "AK2" = 0.7659180977199874
    "AL2" = 0.7521485740564573
    "AM2" = 0.7383790503929272
    "AN2" = 0.724609526729397
    "AO2" = 0.7108400030658669
    "AP2" = 0.6970704794023368
    "AQ2" = 0.6833009557388066
    "AR2" = 0.6695314320752765
    "AS2" = 0.6557619084117463
}
foreach ($key in $values_6.Keys) {
    $newSheet.Range($key).Value = $values_6[$key]
}

# --- Add strategy_id-5013 ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "strategy_id-5013"
$template.Range("A1:AS2").Copy($newSheet.Range("A1:AS2"))

$values_7 = @{
    "U2" = 0.9831702487644064
    "V2" = 0.9663404975288129
    "W2" = 0.9495107462932192
    "X2" = 0.9326809950576255
    "Y2" = 0.915851243822032
    "Z2" = 0.8990214925864384
    "AA2" = 0.8821917413508448
    "AB2" = 0.8653619901152512
    "AC2" = 0.8485322388796576
    "AD2" = 0.831702487644064
    "AE2" = 0.8148727364084704
    "AF2" = 0.7980429851728768
    "AG2" = 0.7812132339372831
    "AH2" = 0.7643834827016895
    "AI2" = 0.7475537314660958
    "AJ2" = 0.7307239802305022
    "AK2" = 0.7138942289949086
    "AL2" = 0.6970644777593151
    "AM2" = 0.6802347265237214
    "AN2" = 0.6634049752881279
    "AO2" = 0.6465752240525343
    "AP2" = 0.6297454728169407
    "AQ2" = 0.6129157215813471
    "AR2" = 0.5960859703457535
    "AS2" = 0.5792562191101598
}
foreach ($key in $values_7.Keys) {
    $newSheet.Range($key).Value = $values_7[$key]
}

# --- Add strategy_id-5014 ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "strategy_id-5014"
$template.Range("A1:AS2").Copy($newSheet.Range("A1:AS2"))

$values_8 = @{
    "U2" = 0.9888585559026216
    "V2" = 0.9777171118052433
    "W2" = 0.966575667707865
    "X2" = 0.9554342236104866
    "Y2" = 0.9442927795131083
    "Z2" = 0.9331513354157299
    "AA2" = 0.9220098913183516
    "AB2" = 0.9108684472209733
    "AC2" = 0.899727003123595
    "AD2" = 0.8885855590262166
    "AE2" = 0.8774441149288383
    "AF2" = 0.86630267083146
    "AG2" = 0.8551612267340816
    "AH2" = 0.8440197826367033
    "AI2" = 0.8328783385393249
    "AJ2" = 0.8217368944419465
    "AK2" = 0.8105954503445683
    "AL2" = 0.79945400624719
    "AM2" = 0.7883125621498116
    "AN2" = 0.7771711180524332
    "AO2" = 0.7660296739550549
    "AP2" = 0.7548882298576767
    "AQ2" = 0.7437467857602983
    "AR2" = 0.7326053416629199
    "AS2" = 0.7214638975655416
}
foreach ($key in $values_8.Keys) {
    $newSheet.Range($key).Value = $values_8[$key]
}

# --- Add strategy_id-5015 ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "strategy_id-5015"
$template.Range("A1:AS2").Copy($newSheet.Range("A1:AS2"))

$values_9 = @{
    "U2" = 0.9813124864730004
    "V2" = 0.9626249729460008
    "W2" = 0.9439374594190011
    "X2" = 0.9252499458920015
    "Y2" = 0.906562432365002
    "Z2" = 0.8878749188380024
    "AA2" = 0.8691874053110027
    "AB2" = 0.8504998917840031
    "AC2" = 0.8318123782570036
    "AD2" = 0.813124864730004
    "AE2" = 0.7944373512030043
    "AF2" = 0.7757498376760047
    "AG2" = 0.7570623241490051
    "AH2" = 0.7383748106220054
    "AI2" = 0.7196872970950059
    "AJ2" = 0.7009997835680062
    "AK2" = 0.6823122700410066
    "AL2" = 0.663624756514007
    "AM2" = 0.6449372429870075
    "AN2" = 0.6262497294600078
    "AO2" = 0.6075622159330083
    "AP2" = 0.5888747024060086
    "AQ2" = 0.5701871888790091
    "AR2" = 0.5514996753520094
    "AS2" = 0.5328121618250098
}
foreach ($key in $values_9.Keys) {
    $newSheet.Range($key).Value = $values_9[$key]
}
